$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2103.625
$ws.Range("I2").Value = 2371.6667
$ws.Range("K2").Value = 2371.6667
$ws.Range("M2").Value = -2258.6667
$ws.Range("H15").Value = 885.7595
$ws.Range("I15").Value = 885.7595
$ws.Range("K15").Value = 2657.2785
$ws.Range("M15").Value = -2488.2785
$ws.Range("H43").Value = 19956.25
$ws.Range("J43").Value = 22621.875
$ws.Range("L43").Value = 22621.875
$ws.Range("N43").Value = -22759.875
$ws.Range("H62").Value = 6991
$ws.Range("H64").Value = 3417.8635
$ws.Range("J64").Value = 5099.25
$ws.Range("L64").Value = 5099.25
$ws.Range("N64").Value = -5595.25
$ws.Range("H65").Value = 6991
$ws.Range("H67").Value = 3417.8635
$ws.Range("J67").Value = 5099.25
$ws.Range("L67").Value = 5099.25
$ws.Range("N67").Value = -6815.25
$ws.Range("H106").Value = 9379.4375
$ws.Range("I106").Value = 2929.3076
$ws.Range("K106").Value = 2929.3076
$ws.Range("M106").Value = -2298.3076
$ws.Range("H125").Value = 2249.3044
$ws.Range("J125").Value = 3944.6
$ws.Range("L125").Value = 35501.4
$ws.Range("N125").Value = -40421.4
$ws.Range("H138").Value = 6360.7754
$ws.Range("J138").Value = 6146.804
$ws.Range("L138").Value = 18440.412
$ws.Range("N138").Value = -28720.412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11694.726
$ws.Range("I32").Value = 10437.082
$ws.Range("K32").Value = 10437.082
$ws.Range("M32").Value = -10150.082
$ws.Range("H61").Value = 1574.7273
$ws.Range("I61").Value = 1574.7273
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1574.7273
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1362.7273
$ws.Range("N61").ClearContents()
$ws.Range("H97").Value = 6906.9375
$ws.Range("I97").Value = 679.4286
$ws.Range("K97").Value = 679.4286
$ws.Range("M97").Value = -183.4286
$ws.Range("H132").Value = 3495.3547
$ws.Range("I132").Value = 3574.3447
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 10723.0341
$ws.Range("L132").Value = 7050
$ws.Range("M132").Value = -8193.034100000001
$ws.Range("N132").Value = -12110
$ws.Range("H136").Value = 1574.7273
$ws.Range("I136").Value = 1574.7273
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4724.1819
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2174.1819
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5572.75
$ws.Range("J86").Value = 4719.4
$ws.Range("L86").Value = 4719.4
$ws.Range("N86").Value = -6965.4
$ws.Range("H89").Value = 5572.75
$ws.Range("J89").Value = 4719.4
$ws.Range("L89").Value = 23597
$ws.Range("N89").Value = -34829
$ws.Range("H94").Value = 52636116
$ws.Range("I94").Value = 62500884
$ws.Range("J94").Value = 24003.334
$ws.Range("K94").Value = 62500884
$ws.Range("L94").Value = 24003.334
$ws.Range("M94").Value = -62500433
$ws.Range("N94").Value = -24905.334
$ws.Range("H105").Value = 2466.6924
$ws.Range("I105").Value = 2380.5833
$ws.Range("K105").Value = 2380.5833
$ws.Range("M105").Value = -633.5832999999998
$ws.Range("H107").Value = 68987.664
$ws.Range("I107").Value = 144179.86
$ws.Range("J107").Value = 3194.5
$ws.Range("K107").Value = 144179.86
$ws.Range("L107").Value = 3194.5
$ws.Range("M107").Value = -142259.86
$ws.Range("N107").Value = -7034.5
$ws.Range("H112").Value = 44995
$ws.Range("J112").Value = 44995
$ws.Range("L112").Value = 44995
$ws.Range("N112").Value = -47949
$ws.Range("H134").Value = 2005.3043
$ws.Range("I134").Value = 1831.15
$ws.Range("K134").Value = 5493.450000000001
$ws.Range("M134").Value = -2958.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 250001000
$ws.Range("I7").Value = 1000000000
$ws.Range("J7").Value = 1332.6666
$ws.Range("K7").Value = 1000000000
$ws.Range("L7").Value = 1332.6666
$ws.Range("M7").Value = -999999887
$ws.Range("N7").Value = -1558.6666
$ws.Range("H50").Value = 14999.75
$ws.Range("J50").Value = 14999.75
$ws.Range("L50").Value = 14999.75
$ws.Range("N50").Value = -16249.75
$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16472
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290
$ws.Range("H60").Value = 14999.833
$ws.Range("J60").Value = 14999.833
$ws.Range("L60").Value = 14999.833
$ws.Range("N60").Value = -16021.833
$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 15000
$ws.Range("L61").Value = 15000
$ws.Range("N61").Value = -15696
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 46666.668
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19126
$ws.Range("H77").Value = 46666.668
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55632

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 7817.3335
$ws.Range("I58").Value = 4202.5
$ws.Range("J58").Value = 9624.75
$ws.Range("K58").Value = 12607.5
$ws.Range("L58").Value = 28874.25
$ws.Range("M58").Value = -12479.5
$ws.Range("N58").Value = -29130.25
$ws.Range("H98").Value = 740.5
$ws.Range("I98").Value = 603.3333
$ws.Range("K98").Value = 1809.9999
$ws.Range("M98").Value = -311.9999
$ws.Range("H107").Value = 622.73846
$ws.Range("I107").Value = 373.8
$ws.Range("J107").Value = 697.42
$ws.Range("K107").Value = 1121.4
$ws.Range("L107").Value = 2092.26
$ws.Range("M107").Value = 798.5999999999999
$ws.Range("N107").Value = -5932.26
$ws.Range("H122").Value = 1166.3334
$ws.Range("I122").Value = 1166.3334
$ws.Range("K122").Value = 10497.0006
$ws.Range("M122").Value = -8047.000599999999
$ws.Range("H128").Value = 501771.66
$ws.Range("I128").Value = 501771.66
$ws.Range("K128").Value = 1505314.98
$ws.Range("M128").Value = -1500334.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2312.4285
$ws.Range("J113").Value = 2377.1428
$ws.Range("L113").Value = 2377.1428
$ws.Range("N113").Value = -6717.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H20").Value = 12000
$ws.Range("J20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12452
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H132").Value = 6497.8
$ws.Range("I132").Value = 4749.5
$ws.Range("K132").Value = 14248.5
$ws.Range("M132").Value = -11718.5
$ws.Range("H136").Value = 7872.543
$ws.Range("I136").Value = 7665.5356
$ws.Range("J136").Value = 8700.571
$ws.Range("K136").Value = 22996.6068
$ws.Range("L136").Value = 26101.713
$ws.Range("M136").Value = -20446.6068
$ws.Range("N136").Value = -31201.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50466.5
$ws.Range("I2").Value = 50466.5
$ws.Range("K2").Value = 50466.5
$ws.Range("M2").Value = -50354.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 619.75
$ws.Range("I113").Value = 626.3333
$ws.Range("K113").Value = 1878.9999
$ws.Range("M113").Value = 291.0001
$ws.Range("H122").Value = 1539.091
$ws.Range("I122").Value = 1103.3334
$ws.Range("K122").Value = 3310.0002
$ws.Range("M122").Value = -860.0001999999999
$ws.Range("H126").Value = 25035.4
$ws.Range("I126").Value = 25035.4
$ws.Range("K126").Value = 75106.20000000001
$ws.Range("M126").Value = -72636.20000000001
$ws.Range("H132").Value = 7748.5625
$ws.Range("I132").Value = 8387.799999999999
$ws.Range("K132").Value = 25163.4
$ws.Range("M132").Value = -22633.4
$ws.Range("H136").Value = 1283
$ws.Range("I136").Value = 1286.7142
$ws.Range("K136").Value = 3860.1426
$ws.Range("M136").Value = -1310.1426
